$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header / data cells with new text (mirrors shared-string content changes)
$ws.Range("C1").Value = "name"
$ws.Range("C2").Value = "Badlands"
$ws.Range("B2").Value = "A"
$ws.Range("C3").Value = "jump"
$ws.Range("C4").Value = "barrin, master wizard"
$ws.Range("C5").Value = "mox jet"
$ws.Range("C6").Value = "volcanic island"
$ws.Range("C7").Value = "island"
$ws.Range("C8").Value = "watery grave"
$ws.Range("D8").Value = "nm"
$ws.Range("C9").Value = "Power artifact"

# Update the active selection on Sheet1
$ws.Range("C9").Select()
